# Insert a new weekly price record as row 305 ("Fruta / hortaliza, semanal").
# This pushes the former rows 305-377 down to become rows 306-378,
# and extends the used range from A1:T377 to A1:T378.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 305, shifting rows 305:377 down to 306:378.
$ws.Rows("305:305").Insert()

# Populate the newly inserted row 305 with the new record's data.
$ws.Cells.Item(305, 1).Value  = 3
$ws.Cells.Item(305, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(305, 3).Value  = "Coquimbo"
$ws.Cells.Item(305, 4).Value  = 45204
$ws.Cells.Item(305, 5).Value  = 5
$ws.Cells.Item(305, 6).Value  = "Fruta"
$ws.Cells.Item(305, 7).Value  = 100101
$ws.Cells.Item(305, 8).Value  = "Berries"
$ws.Cells.Item(305, 9).Value  = 100101001
$ws.Cells.Item(305, 10).Value = "Arándano (blue)"
$ws.Cells.Item(305, 11).Value = "Sin especificar"
$ws.Cells.Item(305, 12).Value = "Primera"
$ws.Cells.Item(305, 13).Value = 56
$ws.Cells.Item(305, 14).Value = 15000
$ws.Cells.Item(305, 15).Value = 15000
$ws.Cells.Item(305, 16).Value = 15000
$ws.Cells.Item(305, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(305, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(305, 19).Value = 7500
$ws.Cells.Item(305, 20).Value = 2
